# Quarterly update: add 01-04-2021 data row and revise the 01-01-2021 row
# ("Actualización desde MV -datos-").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Revise existing row 54 (Serie = 01-01-2021) with updated figures ---
$ws.Cells.Item(54, 2).Value = 110.4   # B54
$ws.Cells.Item(54, 5).Value = 92.3    # E54
$ws.Cells.Item(54, 6).Value = 63.2    # F54
$ws.Cells.Item(54, 7).Value = 62.8    # G54
$ws.Cells.Item(54, 8).Value = 60.5    # H54
$ws.Cells.Item(54, 9).Value = 66.3    # I54
$ws.Cells.Item(54, 10).Value = 109.9  # J54
$ws.Cells.Item(54, 12).Value = 97.8   # L54

# --- Append the new quarter row 55 (Serie = 01-04-2021) ---
# Leading apostrophe forces the date-like label to be stored as literal
# text (like the other "Serie" labels) instead of being auto-converted to
# a date serial; resetting the Style back to Normal drops the transient
# quote-prefix formatting so the cell keeps the sheet's default style.
$ws.Cells.Item(55, 1).Value = "'01-04-2021"
$ws.Cells.Item(55, 1).Style = "Normal"
$ws.Cells.Item(55, 2).Value = 116.2
$ws.Cells.Item(55, 3).Value = 116.6
$ws.Cells.Item(55, 4).Value = 115.8
$ws.Cells.Item(55, 5).Value = 99.40000000000001
$ws.Cells.Item(55, 6).Value = 68.40000000000001
$ws.Cells.Item(55, 7).Value = 63.8
$ws.Cells.Item(55, 8).Value = 54.9
$ws.Cells.Item(55, 9).Value = 87.2
$ws.Cells.Item(55, 10).Value = 118.3
$ws.Cells.Item(55, 11).Value = 101
$ws.Cells.Item(55, 12).Value = 104.2
